$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in the newly journaled work entries (rows 24-27) ---
$ws.Range("A24").Value = "05/03/2018"
$ws.Range("B24").Value = "Débugging Pencil et Eraser"
$ws.Range("C24").Value = 6

$ws.Range("A25").Value = "05/04/2018"
$ws.Range("B25").Value = "fin débugging Pencil et Eraser"
$ws.Range("C25").Value = 6

$ws.Range("A26").Value = "05/05/2018"
$ws.Range("B26").Value = "Ajout de l'outil Pipette et zoom"
$ws.Range("C26").Value = 3

$ws.Range("A27").Value = "05/06/2018"
$ws.Range("B27").Value = "Suite ajout de l'outil zoom"
$ws.Range("C27").Value = 3

# --- Update the active selection to reflect where the author last worked ---
$ws.Range("B29").Select()

$wb.Save()
